$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Center-align the "Immune Status" column (3rd column) in every row,
# including the header row.
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $cell = $t.Cell($r, 3)
    $cell.Range.Paragraphs.Item(1).Alignment = 1
}

# Rows 11 and 12 of the "Immune Status" column also get their text
# upgraded from "uninfected" to "Uninfected controls". Set the text
# directly on a range trimmed of the end-of-cell marker so the edit
# stays scoped to this single cell (Find/Replace on this runtime
# searches the whole document rather than the supplied range).
$row11Range = $t.Cell(11, 3).Range
$row11Range.End = $row11Range.End - 1
$row11Range.Text = "Uninfected controls"

$row12Range = $t.Cell(12, 3).Range
$row12Range.End = $row12Range.End - 1
$row12Range.Text = "Uninfected controls"
